$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll/select to mirror the manual edit (select entire row 35, then delete it)
$ws.Rows.Item(35).Select()
$ws.Rows.Item(35).EntireRow.Delete()
